$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $s = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $s
}

Set-TextValue 'D2' '31.273.38'
Set-TextValue 'E2' '  +3.01%  '
Set-TextValue 'D3' '1.983.30'
Set-TextValue 'E3' '  +5.99%  '
Set-TextValue 'D4' '0.9983'
Set-TextValue 'E4' '  -0.18%  '
Set-TextValue 'D5' '0.7985'
Set-TextValue 'E5' '  +69.79%  '
Set-TextValue 'D6' '253.29'
Set-TextValue 'E6' '  +3.88%  '
Set-TextValue 'D7' '0.9988'
Set-TextValue 'E7' '  -0.14%  '
Set-TextValue 'D8' '0.3386'
Set-TextValue 'E8' '  +17.84%  '
Set-TextValue 'D9' '25.68'
Set-TextValue 'E9' '  +16.72%  '
Set-TextValue 'D10' '0.06938'
Set-TextValue 'E10' '  +7.59%  '
Set-TextValue 'D11' '0.8328'
Set-TextValue 'E11' '  +14.96%  '
Set-TextValue 'E12' '  +4.33%  '
Set-TextValue 'D13' '1.985.54'
Set-TextValue 'E13' '  +6.14%  '
Set-TextValue 'D14' '100.02'
Set-TextValue 'E14' '  +4.25%  '
Set-TextValue 'D15' '5.473'
Set-TextValue 'E15' '  +6.95%  '
Set-TextValue 'D16' '273.51'
Set-TextValue 'E16' '  -1.94%  '
Set-TextValue 'D17' '31.258.53'
Set-TextValue 'E17' '  +3.01%  '
Set-TextValue 'D18' '13.86'
Set-TextValue 'E18' '  +6.82%  '
Set-TextValue 'D19' '0.000007952'
Set-TextValue 'E19' '  +6.24%  '
Set-TextValue 'D20' '5.726'
Set-TextValue 'E20' '  +9.56%  '
Set-TextValue 'D21' '2.247.66'
Set-TextValue 'E21' '  +6.53%  '
Set-TextValue 'D22' '0.9993'
Set-TextValue 'D23' '0.9982'
Set-TextValue 'E23' '  -0.18%  '
Set-TextValue 'D24' '6.925'
Set-TextValue 'E24' '  +11.23%  '
Set-TextValue 'D25' '9.684'
Set-TextValue 'E25' '  +6.97%  '
Set-TextValue 'D26' '163.98'
Set-TextValue 'E26' '  +0.52%  '
Set-TextValue 'E27' '  +52.58%  '
Set-TextValue 'D28' '19.81'
Set-TextValue 'E28' '  +6.33%  '
Set-TextValue 'D29' '2.180'
Set-TextValue 'E29' '  +16.38%  '
Set-TextValue 'D30' '1.566'
Set-TextValue 'E30' '  +6.02%  '
Set-TextValue 'B31' 'Filecoin'
Set-TextValue 'C31' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D31' '4.577'
Set-TextValue 'E31' '  +8.83%  '
Set-TextValue 'B32' 'Toncoin'
Set-TextValue 'C32' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D32' '1.350'
Set-TextValue 'E32' '  +2.22%  '
Set-TextValue 'D33' '4.337'
Set-TextValue 'E33' '  +5.52%  '
Set-TextValue 'D34' '0.05162'
Set-TextValue 'E34' '  +7.40%  '
Set-TextValue 'D35' '1.211'
Set-TextValue 'E35' '  +8.27%  '
Set-TextValue 'D36' '0.7520'
Set-TextValue 'E36' '  +9.14%  '
Set-TextValue 'D37' '2.761'
Set-TextValue 'E37' '  +1.65%  '
Set-TextValue 'D38' '0.9991'
Set-TextValue 'E38' '  -0.13%  '
Set-TextValue 'D39' '0.02008'
Set-TextValue 'E39' '  +7.06%  '
Set-TextValue 'D40' '2.908'
Set-TextValue 'E40' '  +3.45%  '
Set-TextValue 'D41' '6.601'
Set-TextValue 'E41' '  +6.21%  '
Set-TextValue 'D42' '78.07'
Set-TextValue 'E42' '  +5.06%  '
Set-TextValue 'D43' '0.4658'
Set-TextValue 'E43' '  +10.21%  '
Set-TextValue 'D44' '2.058'
Set-TextValue 'E44' '  +6.71%  '
Set-TextValue 'D45' '0.8534'
Set-TextValue 'E45' '  +2.96%  '
Set-TextValue 'D46' '104.77'
Set-TextValue 'E46' '  +3.86%  '
Set-TextValue 'D47' '0.9994'
Set-TextValue 'E47' '  +0.01%  '
Set-TextValue 'D48' '10.01'
Set-TextValue 'E48' '  +4.37%  '
Set-TextValue 'D49' '7.483'
Set-TextValue 'E49' '  +7.75%  '
Set-TextValue 'D50' '0.4284'
Set-TextValue 'E50' '  +9.36%  '
Set-TextValue 'D51' '36.50'
Set-TextValue 'E51' '  +3.42%  '
